$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new text value, driven by the authoritative diff.
$updates = @{
  '2' = @{ 'D'='61.689.80'; 'E'='  -5.28%  ' }
  '3' = @{ 'D'='2.986.92'; 'E'='  -6.80%  ' }
  '4' = @{ 'D'='1.00'; 'E'='  +0.05%  ' }
  '5' = @{ 'D'='544.37'; 'E'='  -5.33%  ' }
  '6' = @{ 'D'='153.19'; 'E'='  -8.49%  ' }
  '7' = @{ 'E'='  -0.03%  ' }
  '8' = @{ 'D'='0.562'; 'E'='  -5.77%  ' }
  '9' = @{ 'D'='2.988.05'; 'E'='  -6.51%  ' }
  '10' = @{ 'E'='  -6.18%  ' }
  '11' = @{ 'D'='6.19'; 'E'='  -8.30%  ' }
  '12' = @{ 'E'='  -6.87%  ' }
  '13' = @{ 'D'='3.509.39'; 'E'='  -6.78%  ' }
  '14' = @{ 'E'='  -3.67%  ' }
  '15' = @{ 'D'='61.770.03'; 'E'='  -5.11%  ' }
  '16' = @{ 'D'='23.54'; 'E'='  -8.21%  ' }
  '17' = @{ 'D'='2.992.38'; 'E'='  -6.48%  ' }
  '18' = @{ 'E'='  -6.76%  ' }
  '19' = @{ 'D'='5.11'; 'E'='  -4.29%  ' }
  '20' = @{ 'D'='387.97'; 'E'='  -5.86%  ' }
  '21' = @{ 'D'='11.90'; 'E'='  -7.87%  ' }
  '22' = @{ 'D'='6.61'; 'E'='  -8.15%  ' }
  '23' = @{ 'E'='  -0.09%  ' }
  '24' = @{ 'D'='64.97'; 'E'='  -6.47%  ' }
  '25' = @{ 'D'='0.467'; 'E'='  -5.05%  ' }
  '26' = @{ 'D'='0.187'; 'E'='  -7.39%  ' }
  '27' = @{ 'B'='Binance-PegBSC-USD'; 'C'='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; 'D'='1.00'; 'E'='  +0.29%  ' }
  '28' = @{ 'B'='PEPE'; 'C'='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; 'D'='0.0₃0938'; 'E'='  -10.66%  ' }
  '29' = @{ 'B'='InternetComputer(DFINITY)'; 'C'='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; 'D'='8.32'; 'E'='  -6.75%  ' }
  '31' = @{ 'D'='1.71'; 'E'='  -7.77%  ' }
  '32' = @{ 'D'='20.28'; 'E'='  -5.95%  ' }
  '33' = @{ 'D'='158.76'; 'E'='  +1.39%  ' }
  '34' = @{ 'D'='6.00'; 'E'='  -6.39%  ' }
  '35' = @{ 'D'='4.59'; 'E'='  -7.68%  ' }
  '36' = @{ 'E'='  -6.72%  ' }
  '37' = @{ 'E'='  -7.02%  ' }
  '38' = @{ 'E'='  -8.35%  ' }
  '39' = @{ 'D'='2.432.70'; 'E'='  -11.72%  ' }
  '40' = @{ 'E'='  -6.52%  ' }
  '41' = @{ 'D'='37.11'; 'E'='  -5.06%  ' }
  '42' = @{ 'D'='22.23'; 'E'='  -8.63%  ' }
  '43' = @{ 'E'='  -7.62%  ' }
  '44' = @{ 'E'='  -6.36%  ' }
  '45' = @{ 'E'='  -0.07%  ' }
  '46' = @{ 'E'='  -6.75%  ' }
  '47' = @{ 'D'='4.98'; 'E'='  -12.13%  ' }
  '48' = @{ 'D'='0.0957'; 'E'='  -3.66%  ' }
  '49' = @{ 'D'='19.68'; 'E'='  -8.72%  ' }
  '50' = @{ 'E'='  +0.05%  ' }
  '51' = @{ 'D'='264.89'; 'E'='  -10.78%  ' }
}

foreach ($row in $updates.Keys) {
  foreach ($col in $updates[$row].Keys) {
    $cellRef = "$col$row"
    $newVal = $updates[$row][$col]
    # Force text storage so numeric-looking strings (e.g. "1.00", "6.19")
    # are not silently coerced into numbers by Excel, then restore the
    # default "Normal" style so no stray number-format style lingers.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newVal
    $cell.Style = "Normal"
  }
}